# Apply the "measures.xlsx" update to the Federal sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Federal")

# --- Text corrections (shared string renames) ---

# Header row: capitalize "universities..." (this label is shared verbatim by the
# Federal, TI and VS sheets, so every sheet that uses it must be updated too,
# otherwise the old text would remain as an orphaned shared string)
$ws.Range("H1").Value = "Universities and other educational establishments"
$wsTI = $wb.Worksheets.Item("TI")
$wsTI.Range("F1").Value = "Universities and other educational establishments"
$wsVS = $wb.Worksheets.Item("VS")
$wsVS.Range("E1").Value = "Universities and other educational establishments"

# Category 1 row, "Singing allowed" column
$ws.Range("W3").Value = "ban for some people (professionals, amateur, adults)"

# Category 2 row, "Borders" column
$ws.Range("B4").Value = "Borders closed to one direct neighbor country"

# Category 3 row, "Borders" column
$ws.Range("B5").Value = "Borders closed to all direct neighbor countries"

# Category 4 row, "Borders" column: the measure text is removed entirely
$ws.Range("B6").ClearContents()

# --- Data value corrections ---

# 2020-06-15: Borders category changes from 1 to 2
$ws.Range("B19").Value = 2

# 2020-12-21: Borders value is removed
$ws.Range("B29").ClearContents()

# 2020-12-28: Borders value removed; Universities/higher-ed category changes from 3 to 2
$ws.Range("B31").ClearContents()
$ws.Range("H31").Value = 2

# --- New/split rows at the bottom of the table (old row 37 split into new rows 37-38, plus a brand new row 39) ---

# New row 37: 2021-03-01, Cultural/Sport facilities category 2, Singing allowed category 1
$ws.Range("A37").Value = 44256
$ws.Range("A37").NumberFormat = "m/d/yyyy"
$ws.Range("D37").ClearContents()
$ws.Range("M37").ClearContents()
$ws.Range("S37").Value = 2
$ws.Range("T37").Value = 2
$ws.Range("W37").Value = 1

# Row 38 (previously row 37): 2021-03-22, Gatherings category 4, Shops/Markets category 1
$ws.Range("A38").Value = 44277
$ws.Range("A38").NumberFormat = "m/d/yyyy"
$ws.Range("D38").Value = 4
$ws.Range("M38").Value = 1

# New row 39: 2021-04-19
$ws.Range("A39").Value = 44305
$ws.Range("A39").NumberFormat = "m/d/yyyy"
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 3
$ws.Range("H39").Value = 1
$ws.Range("K39").Value = 2
$ws.Range("S39").Value = 1
$ws.Range("T39").Value = 1
$ws.Range("U39").Value = 1
